# Add a new "2022-Q1" sheet (fund holdings detail) right before the "总计"
# summary sheet, and prepend a matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" worksheet, positioned just before "总计".
#    The existing "总计" sheet is renamed to "2022-Q1" and repurposed to
#    hold the new per-fund detail, while a fresh "总计" summary sheet is
#    appended right after it (mirrors how the workbook was edited by
#    hand: rename the old tab, then add a new running-total tab).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Name = "2022-Q1"
$total = $wb.Worksheets.Add($null, $ws)
$total.Name = "总计"

# Header row (bold, centered, thin border - matches the other quarter sheets)
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 2]
}
$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fund rows, in descending order of 持有市值(亿元) - same layout used by the
# existing per-quarter sheets (2021-Q1 .. 2021-Q4).
$funds = @(
    @('005669','前海开源公用事业行业股票','258.16','94.53','6.55','16.9095',3),
    @('001875','前海开源沪港深优势精选灵活配置混合','82.95','93.80','9.36','7.7641',3),
    @('010717','前海开源优质企业6个月持有期混合A','58.52','93.93','9.05','5.2961',3),
    @('001837','前海开源沪港深蓝筹精选灵活配置混合','15.15','94.10','9.47','1.4347',2),
    @('001605','富兰克林国海沪港深成长精选股票','39.43','86.45','2.94','1.1592',5),
    @('001874','前海开源沪港深价值精选灵活配置混合','9.62','94.37','9.31','0.8956',3),
    @('004263','华安沪港深机会灵活配置混合','13.52','94.79','5.75','0.7774',7),
    @('000934','国富大中华精选混合QDII','25.71','83.59','2.87','0.7379',8),
    @('006370','国富大中华精选混合QDII美元','25.71','83.59','2.87','0.7379',8),
    @('012588','南方港股通优势企业混合型证券投资基金A','37.54','71.00','1.81','0.6795',10),
    @('013270','前海开源聚利一年持有混合A','7.55','80.39','8.70','0.6568',2),
    @('009846','富兰克林国海港股通远见价值混合','19.47','86.72','3.31','0.6445',3),
    @('010718','前海开源优质企业6个月持有期混合C','7.00','93.93','9.05','0.6335',3),
    @('006039','富兰克林国海估值优势灵活配置混合','17.25','87.22','3.12','0.5382',3),
    @('008546','南方产业优势两年持有期混合A','24.91','74.23','1.99','0.4957',6),
    @('010751','宝盈优质成长混合A','5.64','92.80','8.73','0.4924',2),
    @('008381','前海开源新兴产业混合','6.18','93.63','5.06','0.3127',9),
    @('013123','汇添富精选核心优势一年持有混合A','6.15','66.61','4.20','0.2583',5),
    @('009152','南方瑞盛三年持有期混合A','11.98','77.05','2.00','0.2396',6),
    @('457001','国富亚洲机会股票 (QDII)','5.93','77.36','2.84','0.1684',6),
    @('005228','汇添富港股通专注成长混合','3.64','80.11','4.50','0.1638',5),
    @('010752','宝盈优质成长混合C','0.78','92.80','8.73','0.0681',2),
    @('009781','南方产业优势两年持有期混合C','2.82','74.23','1.99','0.0561',6),
    @('006768','华安沪港深优选混合','0.84','93.09','5.95','0.0500',6),
    @('241001','华宝海外中国混合(QDII)','0.83','86.89','4.77','0.0396',7),
    @('012589','南方港股通优势企业混合型证券投资基金C','2.05','71.00','1.81','0.0371',10),
    @('013271','前海开源聚利一年持有混合C','0.42','80.39','8.70','0.0365',2),
    @('006923','前海开源沪港深非周期性行业股票A','0.54','93.77','4.95','0.0267',9),
    @('009153','南方瑞盛三年持有期混合C','1.02','77.05','2.00','0.0204',6),
    @('013124','汇添富精选核心优势一年持有混合C','0.30','66.61','4.20','0.0126',5),
    @('006924','前海开源沪港深非周期性行业股票C','0.22','93.77','4.95','0.0109',9),
    @('012315','创金合信港股通成长股票型发起式证券投资基金A','0.19','83.49','4.46','0.0085',10),
    @('012316','创金合信港股通成长股票型发起式证券投资基金C','0.10','83.49','4.46','0.0045',10)
)

$row = 2
foreach ($f in $funds) {
    $ws.Cells.Item($row, 1).Value = $row - 2
    $ws.Cells.Item($row, 1).Font.Bold = $true
    $ws.Cells.Item($row, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($row, 1).VerticalAlignment = -4160
    $ws.Cells.Item($row, 1).Borders.LineStyle = 1

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $f[0]
    $ws.Cells.Item($row, 3).Value = $f[1]

    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $f[2]
    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value = $f[3]
    $ws.Cells.Item($row, 6).NumberFormat = "@"
    $ws.Cells.Item($row, 6).Value = $f[4]
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = $f[5]

    $ws.Cells.Item($row, 8).Value = $f[6]

    $row = $row + 1
}

$ws.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Build the new "总计" sheet: same header as before, plus a new
#    "2022-Q1" row on top of the existing per-quarter summary rows.
# ---------------------------------------------------------------------
$total.Cells.Item(1, 2).Value = "日期"
$total.Cells.Item(1, 3).Value = "持有数量(只)"
$total.Cells.Item(1, 4).Value = "持有市值(亿元)"
$totalHeaderRange = $total.Range("B1:D1")
$totalHeaderRange.Font.Bold = $true
$totalHeaderRange.HorizontalAlignment = -4108
$totalHeaderRange.VerticalAlignment = -4160
$totalHeaderRange.Borders.LineStyle = 1

$summary = @(
    @("2022-Q1", 33, 41.37),
    @("2021-Q4", 102, 77.08),
    @("2021-Q3", 69, 56.05),
    @("2021-Q2", 7, 0.18),
    @("2021-Q1", 6, 0.19)
)

$r = 2
foreach ($s in $summary) {
    $total.Cells.Item($r, 1).Value = $r - 2
    $total.Cells.Item($r, 1).Font.Bold = $true
    $total.Cells.Item($r, 1).HorizontalAlignment = -4108
    $total.Cells.Item($r, 1).VerticalAlignment = -4160
    $total.Cells.Item($r, 1).Borders.LineStyle = 1

    $total.Cells.Item($r, 2).Value = $s[0]
    $total.Cells.Item($r, 3).Value = $s[1]
    $total.Cells.Item($r, 4).Value = $s[2]

    $r = $r + 1
}

$total.Range("A1").Select()
